# Balance.xlsx — fix weights-and-balance analyses
# - Adds two new summary rows on GLOBAL RESULTS ("Max aft Xcg MAC" / "Max forward Xcg MAC")
# - Drops the SFORZA rows from the WING sheet's Xcg/Ycg estimation-method comparisons
# - Recomputes several downstream cg values across sheets

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# GLOBAL RESULTS
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("GLOBAL RESULTS")

$ws.Cells.Item(2,3).Value  = -16.064114310967643
$ws.Cells.Item(3,3).Value  = 19.81783431650001
$ws.Cells.Item(4,3).Value  = -0.5451193836938566

$ws.Cells.Item(6,3).Value  = 81.12294163972295
$ws.Cells.Item(7,3).Value  = 22.807390723286282
$ws.Cells.Item(8,3).Value  = -0.062149507186539554

$ws.Cells.Item(10,3).Value = 81.12294163972295
$ws.Cells.Item(11,3).Value = 22.807390723286282
$ws.Cells.Item(12,3).Value = -0.062149507186539554

$ws.Cells.Item(14,3).Value = 1.2451627130101575
$ws.Cells.Item(15,3).Value = 20.35028238293929
$ws.Cells.Item(16,3).Value = -0.04189664103200362

$ws.Cells.Item(18,3).Value = 9.49696099842701
$ws.Cells.Item(19,3).Value = 20.60411470855467
$ws.Cells.Item(20,3).Value = -0.33896271527218036

$ws.Cells.Item(22,1).Value = "Max aft Xcg MAC"
$ws.Cells.Item(22,2).Value = "%"
$ws.Cells.Item(22,3).Value = -1.918191911578765

$ws.Cells.Item(23,1).Value = "Max forward Xcg MAC"
$ws.Cells.Item(23,2).Value = "%"
$ws.Cells.Item(23,3).Value = 83.87257400485379

# ---------------------------------------------------------------------------
# WING — drop the SFORZA comparison rows, keep only TORENBEEK_1982
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("WING")

$ws.Cells.Item(2,3).Value  = 2.1829837889405006
$ws.Cells.Item(6,3).Value  = 21.272983788940497

# Row 11 used to be SFORZA; it becomes the TORENBEEK_1982 row (was row 12)
$ws.Cells.Item(11,1).Value = "TORENBEEK_1982"
$ws.Cells.Item(11,2).Value = "m"
$ws.Cells.Item(11,3).Value = 2.1829837889405006

# Old row 12 (TORENBEEK_1982) becomes a blank separator
$ws.Cells.Item(12,1).Value = " "
$ws.Cells.Item(12,2).ClearContents()
$ws.Cells.Item(12,3).ClearContents()

# Row 13 stays "Ycg ESTIMATION METHOD COMPARISON" (just a re-pointed shared string)
$ws.Cells.Item(13,1).Value = "Ycg ESTIMATION METHOD COMPARISON"

# Row 14 used to be SFORZA (Ycg); it becomes the TORENBEEK_1982 row (was row 16)
$ws.Cells.Item(14,1).Value = "TORENBEEK_1982"
$ws.Cells.Item(14,2).Value = "m"
$ws.Cells.Item(14,3).Value = 6.183898882356097

# Old rows 15 and 16 are no longer part of the table
$ws.Rows.Item(15).Delete()
$ws.Rows.Item(15).Delete()

# ---------------------------------------------------------------------------
# HORIZONTAL TAIL
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("HORIZONTAL TAIL")
$ws.Cells.Item(2,3).Value  = 1.9104625588199173
$ws.Cells.Item(6,3).Value  = 33.45046255881991
$ws.Cells.Item(11,3).Value = 1.9104625588199173

# ---------------------------------------------------------------------------
# VERTICAL TAIL
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("VERTICAL TAIL")
$ws.Cells.Item(8,3).Value = 3.696399999999999

# ---------------------------------------------------------------------------
# LANDING GEARS
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("LANDING GEARS")
$ws.Cells.Item(2,3).Value = 18.071522634084396

Write-Output "edits applied"
